$wb = $excel.ActiveWorkbook

# --- Video5 sheet: the existing hyperlink/url in A1 is removed (cell becomes empty) ---
$ws5 = $wb.Worksheets.Item("Video5")
$ws5.Range("A1").Hyperlinks.Delete()
$ws5.Range("A1").Clear()
$ws5.Range("H22").Select()

# --- Video4 sheet: its old url is replaced by a new youtu.be link with a hyperlink ---
$ws4 = $wb.Worksheets.Item("Video4")
$cell4 = $ws4.Range("A1")
$cell4.ClearContents()
$ws4.Hyperlinks.Add($cell4, "https://youtu.be/UyXkte02GQQ")
$cell4.Style = "Hyperlink"

# Video4 becomes the active / selected sheet
$ws4.Activate()
$ws4.Select()
